$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" sheet: insert a new "2022-Q3" row at the top of the data
#    (row 2) and shift the rest of the quarters down by one row.
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Row 9 does not exist yet - give it the same look (bold / bordered /
# centered) as the other index cells in column A by copying the format
# from the row above before writing into it.
$wsTotal.Range("A8").Copy()
$wsTotal.Range("A9").PasteSpecial(-4122)

$totalData = @(
    @(0, "2022-Q3", 12, 4.19),
    @(1, "2022-Q2", 79, 22.6),
    @(2, "2022-Q1", 86, 35.82),
    @(3, "2021-Q4", 156, 53.08),
    @(4, "2021-Q3", 146, 50.73),
    @(5, "2021-Q2", 67, 11.32),
    @(6, "2021-Q1", 26, 7.71),
    @(7, "2020-Q4", 43, 14)
)

for ($i = 0; $i -lt $totalData.Length; $i++) {
    $r = $i + 2
    $wsTotal.Range("A$r").Value = $totalData[$i][0]
    $wsTotal.Range("B$r").Value = $totalData[$i][1]
    $wsTotal.Range("C$r").Value = $totalData[$i][2]
    $wsTotal.Range("D$r").Value = $totalData[$i][3]
}

# ------------------------------------------------------------------
# 2) New "2022-Q3" worksheet: duplicate the "2022-Q2" sheet (so the
#    header row / column styling match the other quarter sheets
#    exactly), place it right after "总计", rename it, drop the
#    left-over rows from the source sheet and fill in the real data.
# ------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

$q3data = @(
    @(0, '005927', '创金合信新能源汽车主题股票A', '16.74', '92.78', '8.90', '1.4899', 3),
    @(1, '005928', '创金合信新能源汽车主题股票C', '16.65', '92.78', '8.90', '1.4818', 3),
    @(2, '013160', '创金合信碳中和混合A', '5.11', '93.04', '7.79', '0.3981', 4),
    @(3, '013161', '创金合信碳中和混合C', '3.16', '93.04', '7.79', '0.2462', 4),
    @(4, '005076', '创金合信优选回报灵活配置混合', '2.74', '94.05', '8.13', '0.2228', 8),
    @(5, '011147', '创金合信气候变化责任投资股票C', '1.43', '92.14', '8.58', '0.1227', 8),
    @(6, '011146', '创金合信气候变化责任投资股票A', '1.20', '92.14', '8.58', '0.1030', 8),
    @(7, '580006', '东吴新经济混合A', '1.37', '90.50', '5.34', '0.0732', 10),
    @(8, '012617', '东吴新经济混合C', '0.55', '90.50', '5.34', '0.0294', 10),
    @(9, '005331', '益民优势安享灵活配置混合', '0.82', '89.33', '2.67', '0.0219', 10),
    @(10, '009336', '平安中证500指数增强A', '0.16', '90.93', '1.41', '0.0023', 10),
    @(11, '009337', '平安中证500指数增强C', '0.16', '90.93', '1.41', '0.0023', 10)
)

$lastSourceRow = 80
$lastTargetRow = 1 + $q3data.Length

# Drop the rows of the copied sheet that are beyond what we need.
if ($lastTargetRow -lt $lastSourceRow) {
    $wsQ3.Range("A$($lastTargetRow + 1):H$lastSourceRow").Clear()
}

# Columns B..G hold text (e.g. "16.74", keeping trailing zeros), only
# column A (index) and H (rank) are real numbers.
$wsQ3.Range("B2:G$lastTargetRow").NumberFormat = "@"

for ($i = 0; $i -lt $q3data.Length; $i++) {
    $r = $i + 2
    $row = $q3data[$i]
    $wsQ3.Range("A$r").Value = $row[0]
    $wsQ3.Range("B$r").Value = $row[1]
    $wsQ3.Range("C$r").Value = $row[2]
    $wsQ3.Range("D$r").Value = $row[3]
    $wsQ3.Range("E$r").Value = $row[4]
    $wsQ3.Range("F$r").Value = $row[5]
    $wsQ3.Range("G$r").Value = $row[6]
    $wsQ3.Range("H$r").Value = $row[7]
}

Write-Output "done"
